$wb = $excel.ActiveWorkbook

$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsTasas = $wb.Worksheets.Item("tasas")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.02 = 7594.68 pesos`n✅ 7594.68 pesos = 2.01 = 922.97 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

$wsTasas.Range("N10").Value = 494
$wsTasas.Range("O10").Value = 3751.77
$wsTasas.Range("N12").Value = 3785.95
$wsTasas.Range("O12").Value = 460.1
